# Update New Orleans shard workbook:
#  1. Reorder sheets so "review_info" comes before "hotel_info".
#  2. Insert a new "State" column into "hotel_info" (between Hotel_Name and City)
#     populated with "Louisiana" for the existing data row.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder worksheets: review_info first, hotel_info second ---
$reviewSheet = $wb.Worksheets.Item("review_info")
$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet.Move($hotelSheet)

# --- 2. Insert "State" column into hotel_info sheet ---
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Range("C1").EntireColumn.Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"
